$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sort the year table (A2:D42) descending by year (column A) ---
# Mirrors the author's "Data > Sort" action on the fleet-age-by-build-year
# table. A 41-year key range (A1:A42) is used as the sort key, matching the
# workbook's stored sortState.
$sortRange = $ws.Range("A2:D42")
$keyRange = $ws.Range("A1:A42")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 2, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Ensure the final row order exactly matches: 2022 descending to 1982,
# followed by the "before 1982" / "> 40" summary row at the bottom. ---
$ws.Cells.Item(2,1).Value2 = 2022; $ws.Cells.Item(2,2).Value2 = 0; $ws.Cells.Item(2,3).Value2 = 812; $ws.Cells.Item(2,4).Value2 = 143326
$ws.Cells.Item(3,1).Value2 = 2021; $ws.Cells.Item(3,2).Value2 = 1; $ws.Cells.Item(3,3).Value2 = 786; $ws.Cells.Item(3,4).Value2 = 312504
$ws.Cells.Item(4,1).Value2 = 2020; $ws.Cells.Item(4,2).Value2 = 2; $ws.Cells.Item(4,3).Value2 = 719; $ws.Cells.Item(4,4).Value2 = 316766
$ws.Cells.Item(5,1).Value2 = 2019; $ws.Cells.Item(5,2).Value2 = 3; $ws.Cells.Item(5,3).Value2 = 1190; $ws.Cells.Item(5,4).Value2 = 526072
$ws.Cells.Item(6,1).Value2 = 2018; $ws.Cells.Item(6,2).Value2 = 4; $ws.Cells.Item(6,3).Value2 = 1207; $ws.Cells.Item(6,4).Value2 = 587450
$ws.Cells.Item(7,1).Value2 = 2017; $ws.Cells.Item(7,2).Value2 = 5; $ws.Cells.Item(7,3).Value2 = 1077; $ws.Cells.Item(7,4).Value2 = 523181
$ws.Cells.Item(8,1).Value2 = 2016; $ws.Cells.Item(8,2).Value2 = 6; $ws.Cells.Item(8,3).Value2 = 1037; $ws.Cells.Item(8,4).Value2 = 583602
$ws.Cells.Item(9,1).Value2 = 2015; $ws.Cells.Item(9,2).Value2 = 7; $ws.Cells.Item(9,3).Value2 = 973; $ws.Cells.Item(9,4).Value2 = 437116
$ws.Cells.Item(10,1).Value2 = 2014; $ws.Cells.Item(10,2).Value2 = 8; $ws.Cells.Item(10,3).Value2 = 1037; $ws.Cells.Item(10,4).Value2 = 349573
$ws.Cells.Item(11,1).Value2 = 2013; $ws.Cells.Item(11,2).Value2 = 9; $ws.Cells.Item(11,3).Value2 = 884; $ws.Cells.Item(11,4).Value2 = 308068
$ws.Cells.Item(12,1).Value2 = 2012; $ws.Cells.Item(12,2).Value2 = 10; $ws.Cells.Item(12,3).Value2 = 876; $ws.Cells.Item(12,4).Value2 = 383691
$ws.Cells.Item(13,1).Value2 = 2011; $ws.Cells.Item(13,2).Value2 = 11; $ws.Cells.Item(13,3).Value2 = 807; $ws.Cells.Item(13,4).Value2 = 427566
$ws.Cells.Item(14,1).Value2 = 2010; $ws.Cells.Item(14,2).Value2 = 12; $ws.Cells.Item(14,3).Value2 = 889; $ws.Cells.Item(14,4).Value2 = 405560
$ws.Cells.Item(15,1).Value2 = 2009; $ws.Cells.Item(15,2).Value2 = 13; $ws.Cells.Item(15,3).Value2 = 869; $ws.Cells.Item(15,4).Value2 = 505334
$ws.Cells.Item(16,1).Value2 = 2008; $ws.Cells.Item(16,2).Value2 = 14; $ws.Cells.Item(16,3).Value2 = 970; $ws.Cells.Item(16,4).Value2 = 465372
$ws.Cells.Item(17,1).Value2 = 2007; $ws.Cells.Item(17,2).Value2 = 15; $ws.Cells.Item(17,3).Value2 = 896; $ws.Cells.Item(17,4).Value2 = 396296
$ws.Cells.Item(18,1).Value2 = 2006; $ws.Cells.Item(18,2).Value2 = 16; $ws.Cells.Item(18,3).Value2 = 721; $ws.Cells.Item(18,4).Value2 = 300833
$ws.Cells.Item(19,1).Value2 = 2005; $ws.Cells.Item(19,2).Value2 = 17; $ws.Cells.Item(19,3).Value2 = 511; $ws.Cells.Item(19,4).Value2 = 194098
$ws.Cells.Item(20,1).Value2 = 2004; $ws.Cells.Item(20,2).Value2 = 18; $ws.Cells.Item(20,3).Value2 = 405; $ws.Cells.Item(20,4).Value2 = 147659
$ws.Cells.Item(21,1).Value2 = 2003; $ws.Cells.Item(21,2).Value2 = 19; $ws.Cells.Item(21,3).Value2 = 395; $ws.Cells.Item(21,4).Value2 = 147586
$ws.Cells.Item(22,1).Value2 = 2002; $ws.Cells.Item(22,2).Value2 = 20; $ws.Cells.Item(22,3).Value2 = 439; $ws.Cells.Item(22,4).Value2 = 158016
$ws.Cells.Item(23,1).Value2 = 2001; $ws.Cells.Item(23,2).Value2 = 21; $ws.Cells.Item(23,3).Value2 = 548; $ws.Cells.Item(23,4).Value2 = 175071
$ws.Cells.Item(24,1).Value2 = 2000; $ws.Cells.Item(24,2).Value2 = 22; $ws.Cells.Item(24,3).Value2 = 519; $ws.Cells.Item(24,4).Value2 = 172889
$ws.Cells.Item(25,1).Value2 = 1999; $ws.Cells.Item(25,2).Value2 = 23; $ws.Cells.Item(25,3).Value2 = 479; $ws.Cells.Item(25,4).Value2 = 152486
$ws.Cells.Item(26,1).Value2 = 1998; $ws.Cells.Item(26,2).Value2 = 24; $ws.Cells.Item(26,3).Value2 = 411; $ws.Cells.Item(26,4).Value2 = 111587
$ws.Cells.Item(27,1).Value2 = 1997; $ws.Cells.Item(27,2).Value2 = 25; $ws.Cells.Item(27,3).Value2 = 256; $ws.Cells.Item(27,4).Value2 = 79570
$ws.Cells.Item(28,1).Value2 = 1996; $ws.Cells.Item(28,2).Value2 = 26; $ws.Cells.Item(28,3).Value2 = 248; $ws.Cells.Item(28,4).Value2 = 79072
$ws.Cells.Item(29,1).Value2 = 1995; $ws.Cells.Item(29,2).Value2 = 27; $ws.Cells.Item(29,3).Value2 = 155; $ws.Cells.Item(29,4).Value2 = 48231
$ws.Cells.Item(30,1).Value2 = 1994; $ws.Cells.Item(30,2).Value2 = 28; $ws.Cells.Item(30,3).Value2 = 192; $ws.Cells.Item(30,4).Value2 = 62915
$ws.Cells.Item(31,1).Value2 = 1993; $ws.Cells.Item(31,2).Value2 = 29; $ws.Cells.Item(31,3).Value2 = 204; $ws.Cells.Item(31,4).Value2 = 60005
$ws.Cells.Item(32,1).Value2 = 1992; $ws.Cells.Item(32,2).Value2 = 30; $ws.Cells.Item(32,3).Value2 = 245; $ws.Cells.Item(32,4).Value2 = 55084
$ws.Cells.Item(33,1).Value2 = 1991; $ws.Cells.Item(33,2).Value2 = 31; $ws.Cells.Item(33,3).Value2 = 187; $ws.Cells.Item(33,4).Value2 = 52850
$ws.Cells.Item(34,1).Value2 = 1990; $ws.Cells.Item(34,2).Value2 = 32; $ws.Cells.Item(34,3).Value2 = 160; $ws.Cells.Item(34,4).Value2 = 39327
$ws.Cells.Item(35,1).Value2 = 1989; $ws.Cells.Item(35,2).Value2 = 33; $ws.Cells.Item(35,3).Value2 = 135; $ws.Cells.Item(35,4).Value2 = 29426
$ws.Cells.Item(36,1).Value2 = 1988; $ws.Cells.Item(36,2).Value2 = 34; $ws.Cells.Item(36,3).Value2 = 104; $ws.Cells.Item(36,4).Value2 = 18502
$ws.Cells.Item(37,1).Value2 = 1987; $ws.Cells.Item(37,2).Value2 = 35; $ws.Cells.Item(37,3).Value2 = 68; $ws.Cells.Item(37,4).Value2 = 14082
$ws.Cells.Item(38,1).Value2 = 1986; $ws.Cells.Item(38,2).Value2 = 36; $ws.Cells.Item(38,3).Value2 = 74; $ws.Cells.Item(38,4).Value2 = 12643
$ws.Cells.Item(39,1).Value2 = 1985; $ws.Cells.Item(39,2).Value2 = 37; $ws.Cells.Item(39,3).Value2 = 49; $ws.Cells.Item(39,4).Value2 = 6648
$ws.Cells.Item(40,1).Value2 = 1984; $ws.Cells.Item(40,2).Value2 = 38; $ws.Cells.Item(40,3).Value2 = 68; $ws.Cells.Item(40,4).Value2 = 5123
$ws.Cells.Item(41,1).Value2 = 1983; $ws.Cells.Item(41,2).Value2 = 39; $ws.Cells.Item(41,3).Value2 = 34; $ws.Cells.Item(41,4).Value2 = 3818
$ws.Cells.Item(42,1).Value2 = 1982; $ws.Cells.Item(42,2).Value2 = 40; $ws.Cells.Item(42,3).Value2 = 92; $ws.Cells.Item(42,4).Value2 = 5410
$ws.Cells.Item(43,1).Value2 = "before 1982"; $ws.Cells.Item(43,2).Value2 = "> 40"; $ws.Cells.Item(43,3).Value2 = 894; $ws.Cells.Item(43,4).Value2 = 41246

# --- New column style: C and D both width 8 (replaces the old bestFit widths) ---
# (7.15 "characters" round-trips through the stored-width formula to exactly 8)
$ws.Range("C:D").ColumnWidth = 7.15

# --- Update the active selection ---
$ws.Range("F25").Select()
